$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.916.89'
$ws.Range('E2').Value = '  +1.91%  '
$ws.Range('D3').Value = '1.770.71'
$ws.Range('E3').Value = '  +2.35%  '
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '1.002'
$cell.Style = 'Normal'
$ws.Range('E4').Value = '  -0.16%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '327.77'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +1.67%  '
$ws.Range('E6').Value = '  -0.19%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.4471'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -1.32%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.3558'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +1.13%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.07433'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +1.39%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '42.09'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  +1.47%  '
$ws.Range('E11').Value = '  +2.68%  '
$ws.Range('E12').Value = '  -0.15%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '21.00'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +3.14%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '6.030'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +2.14%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '7.247'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  +3.02%  '
$ws.Range('D16').Value = '1.771.29'
$ws.Range('E16').Value = '  +2.50%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '93.40'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  +2.51%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '0.00001062'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  +1.31%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '0.06436'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +1.71%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '0.9999'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -0.27%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '17.14'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +3.44%  '
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('D23').Value = '27.968.46'
$ws.Range('E23').Value = '  +1.93%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '11.29'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +2.32%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '2.107'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +1.53%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '161.86'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('E27').Value = '  +3.19%  '
$ws.Range('D28').Value = '1.973.16'
$ws.Range('E28').Value = '  +2.26%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '2.173'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +6.41%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '124.98'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  +0.65%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '1.113'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  +6.91%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '0.09214'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +1.25%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '5.660'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +6.09%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '3.678'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  +0.71%  '
$ws.Range('E35').Value = '  +2.51%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.02291'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +1.40%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '0.06181'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +4.01%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.2105'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +2.89%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '0.6334'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +2.00%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '4.973'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  +2.56%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '1.183'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -0.23%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '1.392'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +1.63%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '7.912'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  +2.88%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '13.31'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +2.15%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '3.737'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +1.24%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '0.5893'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +2.03%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '122.74'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +0.87%  '
$ws.Range('E48').Value = '  +2.41%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '0.06909'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +1.27%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '1.139'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +2.92%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '73.02'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  +3.13%  '
